# Update the UserName and Email columns for the Mathan73/74/75 rows to
# Mathan90/91/92 (and the corresponding gmail.com addresses).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# UserName column (C2:C4): Mathan73/74/75 -> Mathan90/91/92
$ws.Range("C2").Value = "Mathan90"
$ws.Range("C3").Value = "Mathan91"
$ws.Range("C4").Value = "Mathan92"

# Email column (D2:D4): mathan73/74/75@gmail.com -> mathan90/91/92@gmail.com
$ws.Range("D2").Value = "mathan90@gmail.com"
$ws.Range("D3").Value = "mathan91@gmail.com"
$ws.Range("D4").Value = "mathan92@gmail.com"
